# The post in row 283 ("「口を開く前に、心を開きなさい」") was removed from the
# source data. Deleting the entire worksheet row shifts every row below it
# up by one (row 284 -> 283, 285 -> 284, ... 458 -> 457) and shrinks the
# used range from A1:C458 to A1:C457, matching the published diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(283).Delete()
